$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: add the new achievement entry (date, name, description)
$ws.Range("B9").Value = (Get-Date -Year 2010 -Month 5 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B9").NumberFormat = "d-mmm"
$ws.Range("C9").Value = "Nguyễn Lê Hoàng Dũng"
$ws.Range("D9").Value = "Tạo cơ sở dữ liệu"

# Widen column C to fit the new name
$ws.Columns.Item(3).ColumnWidth = 25.5
